# Edit 1: Slide 2 - fix wording in second paragraph of the GitHub purpose shape
# ("ואפשר" -> "ומאפשר") and widen the same shape (use case paragraph 4 area).
$p = $ppt.ActivePresentation

$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(4)
$tr2 = $shp2.TextFrame.TextRange
$usecasePara = $tr2.Paragraphs(2,1)
$usecaseRun = $usecasePara.Runs(1,1)
$usecaseRun.Text = "הממשק די פשוט ומאפשר לראות את כל הקבצים והשינויים דרכו."

# Widen the shape so the new wording fits comfortably (width only; position/height unchanged).
$shp2.Width = 5631803 / 12700

# Edit 2: Slide 3 - grammar fix ("הוא" -> "היא") in the Git description paragraph.
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange

$para1 = $tr3.Paragraphs(1,1)
$run2 = $para1.Runs(2,1)
$run2.Text = " היא תוכנה שפותחה על ידי "

# Edit 3: Slide 3 - add a trailing space before the WINDOWS run in the CMD paragraph.
$para4 = $tr3.Paragraphs(4,1)
$run3 = $para4.Runs(3,1)
$run3.Text = " שלה ל "
